# Applies the "Finished project overview in the project plan template" edit:
#   - Title placeholder becomes the real project title (no more highlight).
#   - The "Project Overview" section is rewritten with the real project
#     description; the "Team" bookmark moves to span the (now adjacent)
#     "Team Organization" heading instead of the overview paragraphs.
#   - A table cell's two adjacent, identically-formatted runs are merged
#     into a single run.
$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) { return $p }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Title: "<Project Name>" (yellow highlight) -> "Tool Co-op" (no highlight)
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.End = $titleRange.End - 1   # exclude the paragraph mark
$titleRange.Text = "Tool Co-op"
$titleRange.HighlightColorIndex = 0     # wdNoHighlight - drops <w:highlight> entirely

# ---------------------------------------------------------------------------
# 2. Rewrite the "Project Overview" section (4 paragraphs):
#    - "Project Overview" heading (bookmark "Team" now starts at its end)
#    - the "aims to build..." paragraph (now a single run, new wording)
#    - the italic intro paragraph (new multi-run wording incl. proofErr marks)
#    - the "Team Organization" heading (split run around the relocated
#      "_GoBack" bookmark, bookmark "Team" ends at its end)
# ---------------------------------------------------------------------------
$startPara = Get-ParagraphByText $d "Project Overview"
$endPara   = Get-ParagraphByText $d "Team Organization"
$blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newBlockXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
            <w:r><w:t>Project Overview</w:t></w:r>
            <w:bookmarkStart w:id="1" w:name="Team"/>
          </w:p>
          <w:p>
            <w:r><w:t>This project aims to build a system for renting out tools for a community of users.</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:rPr><w:i/></w:rPr></w:pPr>
            <w:r><w:t>This tool co-op system will automatically handle the management of the tool warehouse who is managed by Joe. The software will handle everything from keeping track of tools in a database to renting out those tools to customers, keeping track of checkout and due dates. The overall project will automate the repetitive tasks that Joe does, but still giving Joe functionality such as checking tools and viewing reports of tools and user accounts.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> The primary priorities for the tool co-op </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>is</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> to provide tools for its community and maintain this business over a long period over time.</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
            <w:r><w:t>Team Organizati</w:t></w:r>
            <w:bookmarkStart w:id="2" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="2"/>
            <w:r><w:t>on</w:t></w:r>
            <w:bookmarkEnd w:id="1"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$blockRange.InsertXML($newBlockXml) | Out-Null

# ---------------------------------------------------------------------------
# 3. Table cell ("Phase 2" row): merge ", " + "Architectural, UI, and DB
#    Design" runs into a single run.
# ---------------------------------------------------------------------------
$targetCell = $null
foreach ($t in $d.Tables) {
    foreach ($cell in $t.Range.Cells) {
        if ($cell.Range.Text.TrimEnd([char]13, [char]7) -eq "Phase 2 - Analysis, Architectural, UI, and DB Design") {
            $targetCell = $cell
        }
    }
}

$cellPara = $targetCell.Range.Paragraphs(1)
$cellRange = $cellPara.Range

$newCellXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Phase 2 - Analysis</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>, Architectural, UI, and DB Design</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$cellRange.InsertXML($newCellXml) | Out-Null
